$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells whose new value is a plain number must be pre-formatted
# as Text so Excel stores them as strings (matching the source inlineStr cells)
# instead of auto-converting them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '42.413.10'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '2.236.27'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '245.02'
$ws.Range('E5').Value = '  -0.86%  '
$ws.Range('E6').Value = '  +0.71%  '
$ws.Range('D7').Value = '74.26'
$ws.Range('E7').Value = '  -3.78%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('D10').Value = '43.07'
$ws.Range('E10').Value = '  +2.20%  '
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('E12').Value = '  +0.02%  '
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').Value = '14.46'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').Value = '0.852'
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').Value = '2.230.06'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('D17').Value = '42.223.19'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '0.0000109'
$ws.Range('E18').Value = '  +10.32%  '
$ws.Range('E19').Value = '  +0.89%  '
$ws.Range('D20').Value = '72.08'
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').Value = '10.21'
$ws.Range('E21').Value = '  +34.85%  '
$ws.Range('D22').Value = '231.42'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').Value = '2.17'
$ws.Range('E23').Value = '  -5.16%  '
$ws.Range('D24').Value = '11.73'
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('D27').Value = '2.30'
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('D28').Value = '2.22'
$ws.Range('E28').Value = '  +6.02%  '
$ws.Range('D29').Value = '166.71'
$ws.Range('E29').Value = '  -2.00%  '
$ws.Range('E30').Value = '  +1.40%  '
$ws.Range('D31').Value = '5.82'
$ws.Range('E31').Value = '  +17.87%  '
$ws.Range('D32').Value = '0.0807'
$ws.Range('E32').Value = '  -3.69%  '
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('D34').Value = '29.69'
$ws.Range('E34').Value = '  -9.25%  '
$ws.Range('E35').Value = '  -0.63%  '
$ws.Range('D36').Value = '4.47'
$ws.Range('E36').Value = '  -0.80%  '
$ws.Range('D37').Value = '0.0308'
$ws.Range('E37').Value = '  +1.66%  '
$ws.Range('D38').Value = '13.25'
$ws.Range('E38').Value = '  -7.94%  '
$ws.Range('E39').Value = '  -0.93%  '
$ws.Range('D40').Value = '5.64'
$ws.Range('E40').Value = '  -4.42%  '
$ws.Range('D41').Value = '63.45'
$ws.Range('E41').Value = '  +3.95%  '
$ws.Range('D42').Value = '0.201'
$ws.Range('E42').Value = '  -1.19%  '
$ws.Range('D43').Value = '8.82'
$ws.Range('E43').Value = '  +1.59%  '
$ws.Range('D44').Value = '105.86'
$ws.Range('E44').Value = '  -6.10%  '
$ws.Range('E45').Value = '  +2.79%  '
$ws.Range('D46').Value = '0.995'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '1.14'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '2.38'
$ws.Range('E48').Value = '  +3.46%  '
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').Value = '2.73'
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('D51').Value = '4.07'
$ws.Range('E51').Value = '  -2.66%  '
